$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-14 Wednesday" "2025-05-15 Thursday"

Replace-Text "970÷2=" "555÷9="
Replace-Text "627÷7=" "666÷3="
Replace-Text "106÷9=" "324÷7="
Replace-Text "939÷4=" "474÷3="
Replace-Text "580÷3=" "461÷5="

Replace-Text "687÷7=" "196÷9="
Replace-Text "244÷3=" "259÷5="
Replace-Text "895÷9=" "145÷4="
Replace-Text "341÷6=" "680÷3="
Replace-Text "978÷6=" "976÷3="

Replace-Text "611÷4=" "889÷5="
Replace-Text "771÷3=" "870÷9="
Replace-Text "815÷3=" "403÷4="
Replace-Text "281÷9=" "158÷6="
Replace-Text "568÷9=" "155÷8="

Replace-Text "557÷9=" "198÷9="
Replace-Text "158÷7=" "877÷3="
Replace-Text "163÷3=" "309÷8="
Replace-Text "251÷7=" "166÷3="
Replace-Text "555÷8=" "964÷7="

Replace-Text "527÷3=" "133÷3="
Replace-Text "988÷6=" "540÷7="
Replace-Text "250÷4=" "945÷8="
Replace-Text "142÷4=" "497÷9="
Replace-Text "879÷6=" "124÷8="
